# Update "想去人数" (F column) values for both the "展览" and "全部类型"
# sheets, which contain duplicate event listings that need to stay in sync.

$wb = $excel.ActiveWorkbook

# row number (key) -> new F-column value
$updates = @{
    2  = 327
    4  = 10426
    8  = 1287
    9  = 7203
    14 = 3211
    16 = 317
    17 = 678
    18 = 125
    22 = 1652
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
